$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Chris Lynn"

# Insert a new column before column A, shifting existing data to the right
$ws.Range("A1").EntireColumn.Insert()

# Fill in the new column A values
$ws.Range("A1").Value = "matchNo"
$ws.Range("A2").Value = "1st"
